$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Near the bottom of the document there are two short paragraphs:
#      count-1 : bold   "Play Age of the Gods Norse: King of Asgard for Free"
#      count   : italic "Read our review of Age of the Gods Norse: King of Asgard..."
#    The bold call-to-action paragraph is removed entirely, and the
#    italic paragraph's text is replaced with a new DALLE image
#    prompt (its italic run formatting is kept).
#    These two edits are done first since nothing below them moves.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count

$lastRange = $d.Paragraphs($count).Range
if ($lastRange.Text -notlike "Read our review of Age of the Gods Norse*") {
    throw "Unexpected document layout: last paragraph is not the italic blurb."
}
# Re-materialize the range from raw offsets: paragraphs whose first
# run is an empty run otherwise only *insert* text at Start instead
# of replacing the whole range when you assign to .Text directly on
# the Paragraph's own Range object.
$lastSpan = $d.Range($lastRange.Start, $lastRange.End)
$newBlurb = 'Prompt for DALLE: Create a cartoon-style feature image for "Age of the Gods Norse: King of Asgard". The image should feature a happy Maya warrior wearing glasses. The background should be inspired by the mythical realm ruled by Odin, with elements of Viking-inspired design and pink clouds at sunset. Make sure to add the game''s title and Playtech''s logo in the image. The image should be eye-catching and colorful to attract online slot game enthusiasts.'
$lastSpan.Text = $newBlurb

$ctaRange = $d.Paragraphs($count - 1).Range
if ($ctaRange.Text -notlike "Play Age of the Gods Norse*") {
    throw "Unexpected document layout: paragraph before last is not the bold CTA."
}
$ctaSpan = $d.Range($ctaRange.Start, $ctaRange.End)
$ctaSpan.Delete()

# ------------------------------------------------------------------
# 2. Insert a new paragraph right after the H1 title at the top of
#    the document with a bold "Meta description" lead-in followed by
#    the (non-bold) description text.
# ------------------------------------------------------------------
$titleRange = $d.Paragraphs(1).Range
if ($titleRange.Text -notlike "Play Age of the Gods Norse*") {
    throw "Unexpected document layout: first paragraph is not the H1 title."
}
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$metaRange = $metaPara.Range

$boldText = "Meta description"
$restText = ": Read our review of Age of the Gods Norse: King of Asgard and play for free. Features include an impressive soundtrack, Progressive Jackpot, and Fury Respins of Odin bonus."
$metaRange.Text = $boldText + $restText

$boldRange = $d.Range($metaRange.Start, $metaRange.Start + $boldText.Length)
$boldRange.Font.Bold = 1
